$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" (strikeouts). Regenerated save_data replaces old
# Strike# counts with actual K counts for each start (rows 2-14).
$newK = @{
    2  = 1
    3  = 5
    4  = 4
    5  = 2
    6  = 5
    7  = 2
    8  = 3
    9  = 3
    10 = 5
    11 = 2
    12 = 4
    13 = 3
    14 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
